$wb = $excel.ActiveWorkbook

# ALC row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 24
$ws.Range("I11").Value = 24
$ws.Range("K11").Value = 24
$ws.Range("M11").Value = 116

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9197.4375
$ws.Range("I32").Value = 8471.625
$ws.Range("J32").Value = 9923.25
$ws.Range("K32").Value = 8471.625
$ws.Range("L32").Value = 9923.25
$ws.Range("M32").Value = -8145.625
$ws.Range("N32").Value = -10575.25

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1217.96
$ws.Range("I40").Value = 1182.3529
$ws.Range("J40").Value = 1293.625
$ws.Range("K40").Value = 1182.3529
$ws.Range("L40").Value = 1293.625
$ws.Range("M40").Value = -1007.3529
$ws.Range("N40").Value = -1643.625

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3279.6
$ws.Range("I76").Value = 3279.6
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3279.6
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2964.6
$ws.Range("N76").ClearContents()

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3279.6
$ws.Range("I79").Value = 3279.6
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3279.6
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2187.6
$ws.Range("N79").ClearContents()

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1629.1666
$ws.Range("I106").Value = 993.75
$ws.Range("J106").Value = 2900
$ws.Range("K106").Value = 993.75
$ws.Range("L106").Value = 2900
$ws.Range("M106").Value = -362.75
$ws.Range("N106").Value = -4162

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1774
$ws.Range("J112").Value = 1998.6666
$ws.Range("L112").Value = 5995.9998
$ws.Range("N112").Value = -8211.9998

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 168043.5
$ws.Range("I127").Value = 1386
$ws.Range("J127").Value = 501358.5
$ws.Range("K127").Value = 4158
$ws.Range("L127").Value = 1504075.5
$ws.Range("M127").Value = 802
$ws.Range("N127").Value = -1513995.5

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2767.81
$ws.Range("I138").Value = 1391.4828
$ws.Range("J138").Value = 3329.972
$ws.Range("K138").Value = 4174.4484
$ws.Range("L138").Value = 9989.916000000001
$ws.Range("M138").Value = 965.5515999999998
$ws.Range("N138").Value = -20269.916

# ARM row 29
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 504.5
$ws.Range("I29").Value = 504.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 504.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -196.5
$ws.Range("N29").ClearContents()

# ARM row 30
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1543.2
$ws.Range("I30").Value = 739
$ws.Range("K30").Value = 739
$ws.Range("M30").Value = -589

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2583.6167
$ws.Range("I32").Value = 2331.3276
$ws.Range("J32").Value = 9900
$ws.Range("K32").Value = 2331.3276
$ws.Range("L32").Value = 9900
$ws.Range("M32").Value = -2044.3276
$ws.Range("N32").Value = -10474

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4457.1665
$ws.Range("I61").Value = 4226
$ws.Range("K61").Value = 4226
$ws.Range("M61").Value = -4014

# ARM row 86
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# ARM row 89
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2554.4783
$ws.Range("I132").Value = 2187.65
$ws.Range("K132").Value = 6562.950000000001
$ws.Range("M132").Value = -4032.950000000001

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4457.1665
$ws.Range("I136").Value = 4226
$ws.Range("K136").Value = 12678
$ws.Range("M136").Value = -10128

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2250
$ws.Range("I86").Value = 2250
$ws.Range("K86").Value = 2250
$ws.Range("M86").Value = -1127

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2250
$ws.Range("I89").Value = 2250
$ws.Range("K89").Value = 11250
$ws.Range("M89").Value = -5634

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1243.0769
$ws.Range("I107").Value = 1212.7916
$ws.Range("K107").Value = 1212.7916
$ws.Range("M107").Value = 707.2084

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2694.4285
$ws.Range("I134").Value = 2393.5
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 7180.5
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -4645.5
$ws.Range("N134").Value = -18570

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 64285.57
$ws.Range("J140").Value = 65999.8
$ws.Range("L140").Value = 65999.8
$ws.Range("N140").Value = -76359.8

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24841.717
$ws.Range("I31").Value = 30454.457
$ws.Range("J31").Value = 6983
$ws.Range("K31").Value = 30454.457
$ws.Range("L31").Value = 6983
$ws.Range("M31").Value = -30159.457
$ws.Range("N31").Value = -7573

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 24841.717
$ws.Range("I34").Value = 30454.457
$ws.Range("J34").Value = 6983
$ws.Range("K34").Value = 30454.457
$ws.Range("L34").Value = 6983
$ws.Range("M34").Value = -30252.457
$ws.Range("N34").Value = -7387

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2609.6
$ws.Range("I58").Value = 2941
$ws.Range("J58").Value = 2112.5
$ws.Range("K58").Value = 2941
$ws.Range("L58").Value = 2112.5
$ws.Range("M58").Value = -2738
$ws.Range("N58").Value = -2518.5

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4225.2856
$ws.Range("I99").Value = 3935.4
$ws.Range("K99").Value = 3935.4
$ws.Range("M99").Value = -2437.4

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 608.15
$ws.Range("I107").Value = 644.3214
$ws.Range("J107").Value = 523.75
$ws.Range("K107").Value = 644.3214
$ws.Range("L107").Value = 523.75
$ws.Range("M107").Value = 1275.6786
$ws.Range("N107").Value = -4363.75

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4225.2856
$ws.Range("I126").Value = 3935.4
$ws.Range("K126").Value = 11806.2
$ws.Range("M126").Value = -9336.200000000001

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 11344.576
$ws.Range("I134").Value = 7624.8623
$ws.Range("K134").Value = 22874.5869
$ws.Range("M134").Value = -20339.5869

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2609.6
$ws.Range("I136").Value = 2941
$ws.Range("J136").Value = 2112.5
$ws.Range("K136").Value = 8823
$ws.Range("L136").Value = 6337.5
$ws.Range("M136").Value = -6273
$ws.Range("N136").Value = -11437.5

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 46625.18
$ws.Range("I131").Value = 400792
$ws.Range("J131").Value = 7273.311
$ws.Range("K131").Value = 1202376
$ws.Range("L131").Value = 21819.933
$ws.Range("M131").Value = -1197336
$ws.Range("N131").Value = -31899.933

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 12254.333
$ws.Range("I136").Value = 1440
$ws.Range("J136").Value = 17661.5
$ws.Range("K136").Value = 4320
$ws.Range("L136").Value = 52984.5
$ws.Range("M136").Value = 780
$ws.Range("N136").Value = -63184.5

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 5330.857
$ws.Range("I140").Value = 5330.857
$ws.Range("K140").Value = 15992.571
$ws.Range("M140").Value = -10812.571

# GSM row 18
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 37080370
$ws.Range("I18").Value = 55570550
$ws.Range("K18").Value = 55570550
$ws.Range("M18").Value = -55570257

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16447.666
$ws.Range("I70").Value = 16167
$ws.Range("J70").Value = 17009
$ws.Range("K70").Value = 16167
$ws.Range("L70").Value = 17009
$ws.Range("M70").Value = -15897
$ws.Range("N70").Value = -17549

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 16447.666
$ws.Range("I73").Value = 16167
$ws.Range("J73").Value = 17009
$ws.Range("K73").Value = 16167
$ws.Range("L73").Value = 17009
$ws.Range("M73").Value = -15231
$ws.Range("N73").Value = -18881

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2437.4443
$ws.Range("I80").Value = 2067.8
$ws.Range("J80").Value = 2899.5
$ws.Range("K80").Value = 2067.8
$ws.Range("L80").Value = 2899.5
$ws.Range("M80").Value = -1069.8
$ws.Range("N80").Value = -4895.5

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2437.4443
$ws.Range("I83").Value = 2067.8
$ws.Range("J83").Value = 2899.5
$ws.Range("K83").Value = 10339
$ws.Range("L83").Value = 14497.5
$ws.Range("M83").Value = -5347
$ws.Range("N83").Value = -24481.5

# GSM row 114
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 29999
$ws.Range("J114").Value = 29999
$ws.Range("L114").Value = 29999
$ws.Range("N114").Value = -38677

# LTW row 38
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 30000
$ws.Range("J38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30820

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6333.6665
$ws.Range("I40").Value = 6000.4
$ws.Range("K40").Value = 6000.4
$ws.Range("M40").Value = -5864.4

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3153.5
$ws.Range("I132").Value = 2926.111
$ws.Range("J132").Value = 5200
$ws.Range("K132").Value = 8778.332999999999
$ws.Range("L132").Value = 15600
$ws.Range("M132").Value = -6248.332999999999
$ws.Range("N132").Value = -20660

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3240.4827
$ws.Range("I136").Value = 3728.5715
$ws.Range("J136").Value = 1959.25
$ws.Range("K136").Value = 11185.7145
$ws.Range("L136").Value = 5877.75
$ws.Range("M136").Value = -8635.7145
$ws.Range("N136").Value = -10977.75

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3333.2856
$ws.Range("I136").Value = 3333.3635
$ws.Range("K136").Value = 10000.0905
$ws.Range("M136").Value = -7450.0905

# WVR row 137
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 41355.8
$ws.Range("J137").Value = 41695
$ws.Range("L137").Value = 41695
$ws.Range("N137").Value = -51895
